$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet to reflect the new "through" date
$ws.Name = "Through 2022-07-24"

# Update the row label for July to reflect the new "through" date
$ws.Range("A8").Value = "July (through 07-24)"

# Update July row (row 8) values for each year column
$ws.Range("B8").Value = 29
$ws.Range("C8").Value = 45
$ws.Range("D8").Value = 54
$ws.Range("E8").Value = 58
$ws.Range("F8").Value = 38
$ws.Range("G8").Value = 105
$ws.Range("H8").Value = 118
$ws.Range("I8").Value = 136

# Update Total row (row 9) values for each year column
$ws.Range("B9").Value = 154
$ws.Range("C9").Value = 293
$ws.Range("D9").Value = 444
$ws.Range("E9").Value = 411
$ws.Range("F9").Value = 289
$ws.Range("G9").Value = 577
$ws.Range("H9").Value = 878
$ws.Range("I9").Value = 942
